# Auto-generated script applying scheduled market-data refresh to Leve profit sheets.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H, I, J, K, L, M, N)
# for specific leve rows across all 8 job sheets. No formulas are used in this workbook;
# all of H:N are plain numeric values, so we set them directly.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H11").Value = 11.6
$ws_ALC.Range("I11").Value = 11.6
$ws_ALC.Range("K11").Value = 11.6
$ws_ALC.Range("M11").Value = 128.4
$ws_ALC.Range("H17").Value = 1558.2354
$ws_ALC.Range("J17").Value = 1558.2354
$ws_ALC.Range("L17").Value = 4674.706200000001
$ws_ALC.Range("N17").Value = -5010.706200000001
$ws_ALC.Range("H29").Value = 700
$ws_ALC.Range("I29").Value = 700
$ws_ALC.Range("K29").Value = 2100
$ws_ALC.Range("M29").Value = -1819
$ws_ALC.Range("H34").Value = 6806.364
$ws_ALC.Range("I34").Value = 6806.364
$ws_ALC.Range("K34").Value = 6806.364
$ws_ALC.Range("M34").Value = -6603.364
$ws_ALC.Range("H36").Value = 6806.364
$ws_ALC.Range("I36").Value = 6806.364
$ws_ALC.Range("K36").Value = 6806.364
$ws_ALC.Range("M36").Value = -6091.364
$ws_ALC.Range("H74").Value = 5897.706
$ws_ALC.Range("I74").Value = 5550.467
$ws_ALC.Range("K74").Value = 5550.467
$ws_ALC.Range("M74").Value = -4614.467
$ws_ALC.Range("H77").Value = 5897.706
$ws_ALC.Range("I77").Value = 5550.467
$ws_ALC.Range("K77").Value = 27752.335
$ws_ALC.Range("M77").Value = -23072.335
$ws_ALC.Range("H100").Value = 3199.5
$ws_ALC.Range("I100").Value = 2049.5
$ws_ALC.Range("J100").Value = 5499.5
$ws_ALC.Range("K100").Value = 2049.5
$ws_ALC.Range("L100").Value = 5499.5
$ws_ALC.Range("M100").Value = -1508.5
$ws_ALC.Range("N100").Value = -6581.5
$ws_ALC.Range("H112").Value = 1620
$ws_ALC.Range("J112").Value = 1620
$ws_ALC.Range("L112").Value = 4860
$ws_ALC.Range("N112").Value = -7076
$ws_ALC.Range("H113").Value = 2260.348
$ws_ALC.Range("I113").Value = 2174.15
$ws_ALC.Range("K113").Value = 2174.15
$ws_ALC.Range("M113").Value = 1079.85
$ws_ALC.Range("H138").Value = 2799.7708
$ws_ALC.Range("J138").Value = 3026
$ws_ALC.Range("L138").Value = 9078
$ws_ALC.Range("N138").Value = -19358

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H34").Value = 34999.5
$ws_ARM.Range("I34").Value = 34999.5
$ws_ARM.Range("K34").Value = 34999.5
$ws_ARM.Range("M34").Value = -34728.5
$ws_ARM.Range("H74").Value = 2438611.2
$ws_ARM.Range("I74").Value = 1545097
$ws_ARM.Range("K74").Value = 1545097
$ws_ARM.Range("M74").Value = -1544223
$ws_ARM.Range("H77").Value = 2438611.2
$ws_ARM.Range("I77").Value = 1545097
$ws_ARM.Range("K77").Value = 7725485
$ws_ARM.Range("M77").Value = -7721117
$ws_ARM.Range("H104").Value = 17884.4
$ws_ARM.Range("J104").Value = 17884.4
$ws_ARM.Range("L104").Value = 17884.4
$ws_ARM.Range("N104").Value = -24872.4
$ws_ARM.Range("H120").Value = 31381
$ws_ARM.Range("I120").Value = 31381
$ws_ARM.Range("K120").Value = 31381
$ws_ARM.Range("M120").Value = -26543
$ws_ARM.Range("H121").Value = 0
$ws_ARM.Range("J121").Value = 0
$ws_ARM.Range("L121").Value = 0
$ws_ARM.Range("N121").ClearContents()

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H86").Value = 2868.3635
$ws_BSM.Range("I86").Value = 2359.5
$ws_BSM.Range("K86").Value = 2359.5
$ws_BSM.Range("M86").Value = -1236.5
$ws_BSM.Range("H89").Value = 2868.3635
$ws_BSM.Range("I89").Value = 2359.5
$ws_BSM.Range("K89").Value = 11797.5
$ws_BSM.Range("M89").Value = -6181.5
$ws_BSM.Range("H94").Value = 6197
$ws_BSM.Range("I94").Value = 1262.6666
$ws_BSM.Range("K94").Value = 1262.6666
$ws_BSM.Range("M94").Value = -811.6666

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H4").Value = 1854675.8
$ws_CRP.Range("I4").Value = 4499.5
$ws_CRP.Range("J4").Value = 2002689.9
$ws_CRP.Range("K4").Value = 4499.5
$ws_CRP.Range("L4").Value = 2002689.9
$ws_CRP.Range("M4").Value = -4387.5
$ws_CRP.Range("N4").Value = -2002913.9
$ws_CRP.Range("H93").Value = 11122.714
$ws_CRP.Range("I93").Value = 6310
$ws_CRP.Range("K93").Value = 6310
$ws_CRP.Range("M93").Value = -4438
$ws_CRP.Range("H125").Value = 74650
$ws_CRP.Range("J125").Value = 74650
$ws_CRP.Range("L125").Value = 74650
$ws_CRP.Range("N125").Value = -79570
$ws_CRP.Range("H132").Value = 4760.654
$ws_CRP.Range("I132").Value = 4499.913
$ws_CRP.Range("J132").Value = 6759.6665
$ws_CRP.Range("K132").Value = 13499.739
$ws_CRP.Range("L132").Value = 20278.9995
$ws_CRP.Range("M132").Value = -10969.739
$ws_CRP.Range("N132").Value = -25338.9995
$ws_CRP.Range("H135").Value = 87770
$ws_CRP.Range("J135").Value = 87770
$ws_CRP.Range("L135").Value = 87770
$ws_CRP.Range("N135").Value = -97910

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H68").Value = 1619.6
$ws_CUL.Range("I68").Value = 1524.5
$ws_CUL.Range("K68").Value = 4573.5
$ws_CUL.Range("M68").Value = -3762.5
$ws_CUL.Range("H71").Value = 1619.6
$ws_CUL.Range("I71").Value = 1524.5
$ws_CUL.Range("K71").Value = 13720.5
$ws_CUL.Range("M71").Value = -9664.5
$ws_CUL.Range("H94").Value = 13997
$ws_CUL.Range("I94").Value = 13997
$ws_CUL.Range("K94").Value = 41991
$ws_CUL.Range("M94").Value = -41315
$ws_CUL.Range("H122").Value = 1160.0714
$ws_CUL.Range("I122").Value = 460.5
$ws_CUL.Range("J122").Value = 2092.8333
$ws_CUL.Range("K122").Value = 4144.5
$ws_CUL.Range("L122").Value = 18835.4997
$ws_CUL.Range("M122").Value = -1694.5
$ws_CUL.Range("N122").Value = -23735.4997
$ws_CUL.Range("H126").Value = 9699.75
$ws_CUL.Range("I126").Value = 9699.75
$ws_CUL.Range("K126").Value = 29099.25
$ws_CUL.Range("M126").Value = -24159.25
$ws_CUL.Range("H131").Value = 650622.0600000001
$ws_CUL.Range("J131").Value = 1516795.9
$ws_CUL.Range("L131").Value = 4550387.699999999
$ws_CUL.Range("N131").Value = -4560467.699999999
$ws_CUL.Range("H134").Value = 8949.5
$ws_CUL.Range("I134").Value = 9646.75
$ws_CUL.Range("K134").Value = 28940.25
$ws_CUL.Range("M134").Value = -23870.25
$ws_CUL.Range("H139").Value = 2537.6428
$ws_CUL.Range("J139").Value = 2545.8333
$ws_CUL.Range("L139").Value = 7637.499899999999
$ws_CUL.Range("N139").Value = -17917.4999

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H51").Value = 0
$ws_GSM.Range("J51").Value = 0
$ws_GSM.Range("L51").Value = 0
$ws_GSM.Range("N51").ClearContents()
$ws_GSM.Range("H126").Value = 5017.154
$ws_GSM.Range("I126").Value = 7134.8335
$ws_GSM.Range("J126").Value = 3202
$ws_GSM.Range("K126").Value = 21404.5005
$ws_GSM.Range("L126").Value = 9606
$ws_GSM.Range("M126").Value = -18934.5005
$ws_GSM.Range("N126").Value = -14546

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 3161.75
$ws_LTW.Range("I40").Value = 3270.5715
$ws_LTW.Range("J40").Value = 2400
$ws_LTW.Range("K40").Value = 3270.5715
$ws_LTW.Range("L40").Value = 2400
$ws_LTW.Range("M40").Value = -3134.5715
$ws_LTW.Range("N40").Value = -2672
$ws_LTW.Range("H46").Value = 1754.8235
$ws_LTW.Range("I46").Value = 2049.5
$ws_LTW.Range("J46").Value = 1736.4062
$ws_LTW.Range("K46").Value = 2049.5
$ws_LTW.Range("L46").Value = 1736.4062
$ws_LTW.Range("M46").Value = -1861.5
$ws_LTW.Range("N46").Value = -2112.4062
$ws_LTW.Range("H50").Value = 16737.572
$ws_LTW.Range("I50").Value = 24498
$ws_LTW.Range("J50").Value = 13633.4
$ws_LTW.Range("K50").Value = 24498
$ws_LTW.Range("L50").Value = 13633.4
$ws_LTW.Range("M50").Value = -23861
$ws_LTW.Range("N50").Value = -14907.4
$ws_LTW.Range("H54").Value = 19742
$ws_LTW.Range("J54").Value = 11484
$ws_LTW.Range("L54").Value = 11484
$ws_LTW.Range("N54").Value = -12772
$ws_LTW.Range("H61").Value = 5241.1875
$ws_LTW.Range("I61").Value = 5057.2666
$ws_LTW.Range("K61").Value = 5057.2666
$ws_LTW.Range("M61").Value = -4855.2666
$ws_LTW.Range("H108").Value = 0
$ws_LTW.Range("J108").Value = 0
$ws_LTW.Range("L108").Value = 0
$ws_LTW.Range("N108").ClearContents()
$ws_LTW.Range("H113").Value = 5241.1875
$ws_LTW.Range("I113").Value = 5057.2666
$ws_LTW.Range("K113").Value = 5057.2666
$ws_LTW.Range("M113").Value = -2887.2666
$ws_LTW.Range("H122").Value = 3494.2856
$ws_LTW.Range("I122").Value = 3221.3333
$ws_LTW.Range("J122").Value = 3568.7273
$ws_LTW.Range("K122").Value = 9663.999899999999
$ws_LTW.Range("L122").Value = 10706.1819
$ws_LTW.Range("M122").Value = -7213.999899999999
$ws_LTW.Range("N122").Value = -15606.1819
$ws_LTW.Range("H132").Value = 5137.067
$ws_LTW.Range("I132").Value = 5269.8184
$ws_LTW.Range("K132").Value = 15809.4552
$ws_LTW.Range("M132").Value = -13279.4552

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H39").Value = 15032.333
$ws_WVR.Range("I39").Value = 14999
$ws_WVR.Range("K39").Value = 14999
$ws_WVR.Range("M39").Value = -14586
$ws_WVR.Range("H54").Value = 19068.445
$ws_WVR.Range("J54").Value = 19577
$ws_WVR.Range("L54").Value = 19577
$ws_WVR.Range("N54").Value = -20617
$ws_WVR.Range("H96").Value = 3097.1333
$ws_WVR.Range("I96").Value = 3828.5557
$ws_WVR.Range("K96").Value = 3828.5557
$ws_WVR.Range("M96").Value = -2455.5557
$ws_WVR.Range("H126").Value = 2460.8667
$ws_WVR.Range("I126").Value = 1867.6666
$ws_WVR.Range("K126").Value = 5602.9998
$ws_WVR.Range("M126").Value = -3132.9998
